$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.443.13"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.575.87"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.83"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.79"
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("D7").Value = "3.574.58"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.16"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "4.184.74"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "3.579.44"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "65.506.25"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.94"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "3.720.60"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.21"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.69"
$ws.Range("E28").Value = "  +32.71%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.08"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.32"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "3.582.37"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.32"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "171.92"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.01"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.19"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.57"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0830"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.831"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.46"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.26"
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.15"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.49"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "2.463.64"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("E51").Value = "  +1.68%  "
